# Applies: "Se crean nuevos casos para cesion de contrato nit a nit"
#  - Updates several ID values on the existing "Semilla 6" sheet
#  - Removes the wrap-text style from columns C/D (rows 9-14) in favour of a
#    plain/underlined style
#  - Duplicates "Semilla 6" into a new "Semilla 3" sheet (placed right after
#    it, and left as the active/selected tab) with its own set of
#    environment values + two extra formatted (but empty) cells
#  - Re-points the hyperlinks on the new sheet to the new environment URLs
#  - Restores the original selection/active-cell on "Semilla 6"

$wb = $excel.ActiveWorkbook
$sheet6 = $wb.Worksheets.Item("Semilla 6")

# ---------------------------------------------------------------------------
# 1) Update the changed id/reference values on "Semilla 6"
# ---------------------------------------------------------------------------
$sheet6.Range("B9").Value  = "585087750"
$sheet6.Range("E9").Value  = "3045984556"
$sheet6.Range("B10").Value = "126601516"
$sheet6.Range("E10").Value = "3052749177"
$sheet6.Range("B11").Value = "256424866"
$sheet6.Range("B12").Value = "681590982"
$sheet6.Range("B13").Value = "884243417"
$sheet6.Range("B14").Value = "884243417"

# Columns C & D (rows 9-14) move from the wrap-text style to a plain
# underlined style
$rng6 = $sheet6.Range("C9:D14")
$rng6.WrapText = $false
$rng6.Font.Underline = $true

# ---------------------------------------------------------------------------
# 2) Duplicate "Semilla 6" -> "Semilla 3" (placed right after it). Excel
#    automatically activates the freshly inserted copy, which matches the
#    target workbook (tabSelected/activeTab moves to the new sheet).
# ---------------------------------------------------------------------------
$sheet6.Copy([System.Reflection.Missing]::Value, $sheet6)
$sheet3 = $wb.Worksheets.Item(2)
$sheet3.Name = "Semilla 3"

# ---------------------------------------------------------------------------
# 3) "Semilla 3" specific environment values
# ---------------------------------------------------------------------------
$sheet3.Range("A2").Value = "http://10.69.60.106:8180/tigo-pos-web/index.jsp"
$sheet3.Range("B2").Value = "http://10.69.60.107:8280/portal/login?initialURI=%2Fportal%2FCRMPortal%2FVenta"
$sheet3.Range("C2").Value = "http://10.69.60.106:8180/tigo-pos-web/wap/windex.wml"

$sheet3.Range("A4").Value = "10.69.60.103"
$sheet3.Range("B4").Value = "DEV11G"

$sheet3.Range("A5").Value = "10.69.60.102"
$sheet3.Range("E5").Value = "10.69.60.105"

$sheet3.Range("A6").Value = "10.69.60.102"

$sheet3.Range("A7").Value = "10.65.32.76"
$sheet3.Range("B7").Value = "SIEBEL02"

# Re-point the hyperlinks that changed target on "Semilla 3" (the copy keeps
# the old ones, so refresh the whole collection with the right URLs)
$sheet3.Range("A1").Hyperlinks.Delete()
$sheet3.Hyperlinks.Add($sheet3.Range("D2"), "http://10.65.45.12:9001/gatewaycbs/BcServicesInt") | Out-Null
$sheet3.Hyperlinks.Add($sheet3.Range("E2"), "http://10.65.45.12:9001/gatewaymgint/GatewayMGWSInt") | Out-Null
$sheet3.Hyperlinks.Add($sheet3.Range("A2"), "http://10.69.60.106:8180/tigo-pos-web/index.jsp") | Out-Null
$sheet3.Hyperlinks.Add($sheet3.Range("B2"), "http://10.69.60.107:8280/portal/login?initialURI=%2Fportal%2FCRMPortal%2FVenta") | Out-Null
$sheet3.Hyperlinks.Add($sheet3.Range("C2"), "http://10.69.60.106:8180/tigo-pos-web/wap/windex.wml") | Out-Null
$sheet3.Hyperlinks.Add($sheet3.Range("I2"), "http://10.69.60.76:8080/PortabilidadServiceEAR-HPNPCommunicationsDelegateEJB/NPCRMWSImpl") | Out-Null

# Two extra formatted (empty) cells further down the new sheet
$sheet3.Range("C16").Font.Underline = $true
$sheet3.Range("E18").Font.Underline = $true

# ---------------------------------------------------------------------------
# 4) Selections matching the target workbook
# ---------------------------------------------------------------------------
$sheet6.Activate()
$sheet6.Range("E15").Select()

$sheet3.Activate()
$sheet3.Range("E14").Select()
